# Updates the "Jogos da Semana" sheet:
#  - Row 2 odds get refreshed values
#  - Row 3 (Moreirense vs Gil Vicente) is removed entirely
#  - Former Row 4 (CA Cerro vs Maldonado) shifts up to become the new Row 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update specific odds values on row 2
$ws.Range("G2").Value = 1.42
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 2
$ws.Range("L2").Value = 8
$ws.Range("Z2").Value = 9
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 8.5
$ws.Range("AE2").Value = 26
$ws.Range("AH2").Value = 41
$ws.Range("AJ2").Value = 101
$ws.Range("AK2").Value = 67
$ws.Range("AO2").Value = 7
$ws.Range("AQ2").Value = 21
$ws.Range("BA2").Value = 251

# Remove row 3 entirely; row 4 (CA Cerro - Maldonado) shifts up to row 3
$ws.Rows("3").Delete()
